$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.961.21"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "2.948.94"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "379.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.31%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.581"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.20%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("E12").Value = "  +1.45%  "
$ws.Range("D13").Value = "3.411.01"
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "12.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +67.95%  "
$ws.Range("D17").Value = "2.941.11"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.997"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("D19").Value = "50.929.96"
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.37%  "
$ws.Range("D22").Value = "0.0₃0950"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +13.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.71%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "25.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.163"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.04%  "
$ws.Range("E31").Value = "  -3.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "33.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0433"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.58%  "
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("E38").Value = "  +4.55%  "
$ws.Range("E39").Value = "  +1.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.01%  "
$ws.Range("E46").Value = "  -1.83%  "
$ws.Range("E47").Value = "  -1.36%  "
$ws.Range("D48").Value = "2.012.16"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("E49").Value = "  -4.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0315"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.88%  "
